$wb = $excel.ActiveWorkbook

# --- tc022 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tc022 = $wb.Worksheets.Add($null, $lastSheet)
$tc022.Name = "tc022"
$tc022.Range("A1").Value = "fieldname"
$tc022.Range("A2").Value = "Desciption"
$tc022.Range("B1").Value = "def_value"
$tc022.Range("B2").Value = "Testing"
$tc022.Columns.Item(1).ColumnWidth = 10.166666666666666
[void]$tc022.Range("E8").Select()

# --- tc023 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tc023 = $wb.Worksheets.Add($null, $lastSheet)
$tc023.Name = "tc023"
$tc023.Range("A1").Value = "fieldname"
$tc023.Range("A2").Value = "New Rq"
$tc023.Range("B1").Value = "value"
$tc023.Range("B2").Value = "New Rq"
$tc023.Range("C1").Value = "text"
$tc023.Range("C2").Value = "Text Box"
[void]$tc023.Range("E3").Select()

# --- tc028 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$tc028 = $wb.Worksheets.Add($null, $lastSheet)
$tc028.Name = "tc028"
$tc028.Range("A1").Value = "fieldname"
$tc028.Range("A2").Value = "Category"
$tc028.Range("B1").Value = "datatype"
$tc028.Range("B2").Value = "Text Box"
[void]$tc028.Range("H5").Select()

[void]$tc028.Activate()
